$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 6315.75
$ws.Range("I62").Value = 7000
$ws.Range("K62").Value = 7000
$ws.Range("M62").Value = -6376

$ws.Range("H65").Value = 6315.75
$ws.Range("I65").Value = 7000
$ws.Range("K65").Value = 35000
$ws.Range("M65").Value = -31880

$ws.Range("H86").Value = 2914.7058
$ws.Range("I86").Value = 4150
$ws.Range("J86").Value = 2050
$ws.Range("K86").Value = 4150
$ws.Range("L86").Value = 2050
$ws.Range("M86").Value = -3027
$ws.Range("N86").Value = -4296

$ws.Range("H89").Value = 2914.7058
$ws.Range("I89").Value = 4150
$ws.Range("J89").Value = 2050
$ws.Range("K89").Value = 20750
$ws.Range("L89").Value = 10250
$ws.Range("M89").Value = -15134
$ws.Range("N89").Value = -21482

$ws.Range("H101").Value = 823.2
$ws.Range("I101").Value = 583
$ws.Range("K101").Value = 1749
$ws.Range("M101").Value = -127

$ws.Range("H106").Value = 12912.818
$ws.Range("I106").Value = 13603.5
$ws.Range("K106").Value = 13603.5
$ws.Range("M106").Value = -12972.5

$ws.Range("H129").Value = 838.04083
$ws.Range("J129").Value = 993.3946999999999
$ws.Range("L129").Value = 2980.1841
$ws.Range("N129").Value = -12980.1841

$ws.Range("H132").Value = 5132029
$ws.Range("I132").Value = 6538940.5
$ws.Range("K132").Value = 19616821.5
$ws.Range("M132").Value = -19614291.5

$ws.Range("H137").Value = 1104.6666
$ws.Range("I137").Value = 860.625
$ws.Range("J137").Value = 1356.5807
$ws.Range("K137").Value = 2581.875
$ws.Range("L137").Value = 4069.7421
$ws.Range("M137").Value = -31.875
$ws.Range("N137").Value = -9169.742099999999

$ws.Range("H138").Value = 1528.42
$ws.Range("I138").Value = 980.3333
$ws.Range("J138").Value = 1545.3711
$ws.Range("K138").Value = 2940.9999
$ws.Range("L138").Value = 4636.1133
$ws.Range("M138").Value = 2199.0001
$ws.Range("N138").Value = -14916.1133

$ws.Range("H141").Value = 788
$ws.Range("I141").Value = 602.2632
$ws.Range("J141").Value = 2552.5
$ws.Range("K141").Value = 1806.7896
$ws.Range("L141").Value = 7657.5
$ws.Range("M141").Value = 3373.2104
$ws.Range("N141").Value = -18017.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 93.333336
$ws.Range("I5").Value = 100
$ws.Range("J5").Value = 80
$ws.Range("K5").Value = 100
$ws.Range("L5").Value = 80
$ws.Range("M5").Value = 12
$ws.Range("N5").Value = -304

$ws.Range("H32").Value = 2530.675
$ws.Range("I32").Value = 2307.1025
$ws.Range("K32").Value = 2307.1025
$ws.Range("M32").Value = -2020.1025

$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("N40").ClearContents()

$ws.Range("H45").Value = 1428.25
$ws.Range("I45").Value = 1414
$ws.Range("J45").Value = 1456.75
$ws.Range("K45").Value = 1414
$ws.Range("L45").Value = 1456.75
$ws.Range("M45").Value = -1037
$ws.Range("N45").Value = -2210.75

$ws.Range("H61").Value = 1132.3137
$ws.Range("I61").Value = 997.6829
$ws.Range("J61").Value = 1684.3
$ws.Range("K61").Value = 997.6829
$ws.Range("L61").Value = 1684.3
$ws.Range("M61").Value = -785.6829
$ws.Range("N61").Value = -2108.3

$ws.Range("H74").Value = 1178.88
$ws.Range("I74").Value = 626.4706
$ws.Range("J74").Value = 2352.75
$ws.Range("K74").Value = 626.4706
$ws.Range("L74").Value = 2352.75
$ws.Range("M74").Value = 247.5294
$ws.Range("N74").Value = -4100.75

$ws.Range("H77").Value = 1178.88
$ws.Range("I77").Value = 626.4706
$ws.Range("J77").Value = 2352.75
$ws.Range("K77").Value = 3132.353
$ws.Range("L77").Value = 11763.75
$ws.Range("M77").Value = 1235.647
$ws.Range("N77").Value = -20499.75

$ws.Range("H97").Value = 482.9524
$ws.Range("I97").Value = 412.73685
$ws.Range("K97").Value = 412.73685
$ws.Range("M97").Value = 83.26315

$ws.Range("H122").Value = 1306.875
$ws.Range("I122").Value = 1145.2858
$ws.Range("K122").Value = 3435.8574
$ws.Range("M122").Value = -985.8574000000003

$ws.Range("H132").Value = 2115.4666
$ws.Range("I132").Value = 1874.9474
$ws.Range("J132").Value = 2530.9092
$ws.Range("K132").Value = 5624.8422
$ws.Range("L132").Value = 7592.7276
$ws.Range("M132").Value = -3094.8422
$ws.Range("N132").Value = -12652.7276

$ws.Range("H136").Value = 1132.3137
$ws.Range("I136").Value = 997.6829
$ws.Range("J136").Value = 1684.3
$ws.Range("K136").Value = 2993.0487
$ws.Range("L136").Value = 5052.9
$ws.Range("M136").Value = -443.0487000000003
$ws.Range("N136").Value = -10152.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 93.333336
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 100
$ws.Range("L4").Value = 80
$ws.Range("M4").Value = 15
$ws.Range("N4").Value = -310

$ws.Range("H53").Value = 30000
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()

$ws.Range("H105").Value = 66668416
$ws.Range("I105").Value = 76924860
$ws.Range("J105").Value = 1505.5
$ws.Range("K105").Value = 76924860
$ws.Range("L105").Value = 1505.5
$ws.Range("M105").Value = -76923113
$ws.Range("N105").Value = -4999.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 227.27272
$ws.Range("I7").Value = 16.25
$ws.Range("J7").Value = 347.85715
$ws.Range("K7").Value = 16.25
$ws.Range("L7").Value = 347.85715
$ws.Range("M7").Value = 96.75
$ws.Range("N7").Value = -573.85715

$ws.Range("H31").Value = 2119.4
$ws.Range("I31").Value = 2294
$ws.Range("J31").Value = 1566.5
$ws.Range("K31").Value = 2294
$ws.Range("L31").Value = 1566.5
$ws.Range("M31").Value = -1999
$ws.Range("N31").Value = -2156.5

$ws.Range("H34").Value = 2119.4
$ws.Range("I34").Value = 2294
$ws.Range("J34").Value = 1566.5
$ws.Range("K34").Value = 2294
$ws.Range("L34").Value = 1566.5
$ws.Range("M34").Value = -2092
$ws.Range("N34").Value = -1970.5

$ws.Range("H58").Value = 800
$ws.Range("I58").Value = 762.4583
$ws.Range("K58").Value = 762.4583
$ws.Range("M58").Value = -559.4583

$ws.Range("H134").Value = 1125.122
$ws.Range("I134").Value = 1244.6774
$ws.Range("J134").Value = 754.5
$ws.Range("K134").Value = 3734.0322
$ws.Range("L134").Value = 2263.5
$ws.Range("M134").Value = -1199.0322
$ws.Range("N134").Value = -7333.5

$ws.Range("H136").Value = 800
$ws.Range("I136").Value = 762.4583
$ws.Range("K136").Value = 2287.3749
$ws.Range("M136").Value = 262.6251000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 2562
$ws.Range("J106").Value = 2562
$ws.Range("L106").Value = 7686
$ws.Range("N106").Value = -9578

$ws.Range("H122").Value = 985.1667
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 19643.803
$ws.Range("I132").Value = 1298.2258
$ws.Range("J132").Value = 42392.32
$ws.Range("K132").Value = 3894.6774
$ws.Range("L132").Value = 127176.96
$ws.Range("M132").Value = -1364.6774
$ws.Range("N132").Value = -132236.96

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2689.0264
$ws.Range("I132").Value = 2919.0344
$ws.Range("K132").Value = 8757.1032
$ws.Range("M132").Value = -6227.1032
